$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C57").Value = "EXTORTION"
$ws.Range("C95").Value = "EXTORTION"
$ws.Range("C159").Value = "EXTORTION"
$ws.Range("C161").Value = "EXTORTION"
$ws.Range("C162").Value = "EXTORTION"
$ws.Range("C164").Value = "EXTORTION"
$ws.Range("C181").Value = "EXTORTION"
$ws.Range("C183").Value = "EXTORTION"
$ws.Range("C184").Value = "EXTORTION"
$ws.Range("C187").Value = "EXTORTION"
$ws.Range("C195").Value = "EXTORTION"
$ws.Range("C197").Value = "EXTORTION"
$ws.Range("C198").Value = "EXTORTION"
$ws.Range("C200").Value = "EXTORTION"
$ws.Range("C204").Value = "EXTORTION"
$ws.Range("C213").Value = "EXTORTION"
$ws.Range("C214").Value = "EXTORTION"
$ws.Range("C218").Value = "EXTORTION"
$ws.Range("C220").Value = "EXTORTION"
$ws.Range("C221").Value = "EXTORTION"
$ws.Range("C227").Value = "EXTORTION"
$ws.Range("C229").Value = "EXTORTION"
$ws.Range("C230").Value = "EXTORTION"
$ws.Range("C231").Value = "EXTORTION"
$ws.Range("C324").Value = "EXTORTION"
$ws.Range("C471").Value = "EXTORTION"
$ws.Range("C569").Value = "EXTORTION"
$ws.Range("C740").Value = "EXTORTION"
$ws.Range("B754").Value = "EXTORTION"
$ws.Range("C756").Value = "NO EXTORTION"
$ws.Range("C765").Value = "NO EXTORTION"
$ws.Range("C784").Value = "NO EXTORTION"
$ws.Range("C796").Value = "NO EXTORTION"
$ws.Range("C800").Value = "EXTORTION"
$ws.Range("C807").Value = "EXTORTION"
$ws.Range("C812").Value = "EXTORTION"
$ws.Range("C816").Value = "NO EXTORTION"
$ws.Range("C819").Value = "NO EXTORTION"
$ws.Range("C822").Value = "NO EXTORTION"
$ws.Range("C844").Value = "EXTORTION"
$ws.Range("C863").Value = "NO EXTORTION"
$ws.Range("C871").Value = "EXTORTION"
$ws.Range("C894").Value = "EXTORTION"
$ws.Range("C895").Value = "NO EXTORTION"
$ws.Range("C908").Value = "NO EXTORTION"
$ws.Range("C919").Value = "NO EXTORTION"
$ws.Range("C925").Value = "NO EXTORTION"
$ws.Range("C929").Value = "EXTORTION"
$ws.Range("C947").Value = "EXTORTION"
$ws.Range("C950").Value = "EXTORTION"
$ws.Range("C953").Value = "EXTORTION"
$ws.Range("C973").Value = "EXTORTION"
$ws.Range("C974").Value = "EXTORTION"
$ws.Range("C975").Value = "EXTORTION"
$ws.Range("C982").Value = "EXTORTION"
